$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'48.262.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.91%  "
$ws.Range("D3").Value = "'2.529.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'323.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").Value = "'109.18"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("E7").Value = "  +0.59%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.559"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.24%  "
$ws.Range("D10").Value = "'40.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.64%  "
$ws.Range("D11").Value = "'20.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +10.34%  "
$ws.Range("D12").Value = "'0.0822"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("D14").Value = "'7.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.32%  "
$ws.Range("D15").Value = "'2.925.05"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").Value = "'2.523.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").Value = "'0.860"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").Value = "'48.166.27"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.86%  "
$ws.Range("D19").Value = "'13.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.63%  "
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").Value = "'0.0₃0947"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("D23").Value = "'72.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.82%  "
$ws.Range("D24").Value = "'270.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.01%  "
$ws.Range("D25").Value = "'2.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("D26").Value = "'26.29"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").Value = "'10.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.14%  "
$ws.Range("D29").Value = "'0.146"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.61%  "
$ws.Range("D30").Value = "'35.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("E31").Value = "  -8.78%  "
$ws.Range("D32").Value = "'49.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.12%  "
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Value = "'0.0791"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("E38").Value = "  +1.23%  "
$ws.Range("D39").Value = "'3.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("D41").Value = "'22.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.78%  "
$ws.Range("D42").Value = "'2.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.13%  "
$ws.Range("D43").Value = "'117.85"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.77%  "
$ws.Range("D44").Value = "'0.0299"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.27%  "
$ws.Range("D45").Value = "'2.010.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.99%  "
$ws.Range("D46").Value = "'3.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.67%  "
$ws.Range("D47").Value = "'1.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.21%  "
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("D49").Value = "'9.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("E50").Value = "  +0.62%  "
$ws.Range("E51").Value = "  +3.06%  "
